$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2..19 (date serial, y_0 year, y_0_forecast, y_1 year, y_1_forecast)
$targetRows = @(
    @{Row=2; A=39400; B=2007; C=1.75539628881467; D=2008; E=0.2337905658324813},
    @{Row=3; A=39765; B=2008; C=2.213911448916162; D=2009; E=3.386383090739953},
    @{Row=4; A=40130; B=2009; C=2.533533936850563; D=2010; E=0.984293482975751},
    @{Row=5; A=40494; B=2010; C=2.088987486264915; D=2011; E=3.612753212925401},
    @{Row=6; A=40862; B=2011; C=1.212544822741002; D=2012; E=2.158838189283174},
    @{Row=7; A=41228; B=2012; C=1.196776590518644; D=2013; E=1.194058515117336},
    @{Row=8; A=41592; B=2013; C=0.4712609263772594; D=2014; E=1.409662779709819},
    @{Row=9; A=41957; B=2014; C=0.8783377572271434; D=2015; E=2.372074663906587},
    @{Row=10; A=42321; B=2015; C=2.29066283401107; D=2016; E=4.595879021798321},
    @{Row=11; A=42689; B=2016; C=4.109890522944348; D=2017; E=4.034919509273061},
    @{Row=12; A=43053; B=2017; C=1.336316831462692; D=2018; E=0.02883756256675252},
    @{Row=13; A=43418; B=2018; C=1.197912858979611; D=2019; E=0.9262553939922924},
    @{Row=14; A=43783; B=2019; C=1.727537197898665; D=2020; E=2.928189816005666},
    @{Row=15; A=44159; B=2020; C=3.647228437274408; D=2021; E=3.673004547855219},
    @{Row=16; A=44525; B=2021; C=2.777797690741424; D=2022; E=1.579011422502852},
    @{Row=17; A=44890; B=2022; C=0.6994919452575576; D=2023; E=-2.087978868409623},
    @{Row=18; A=45254; B=2023; C=-1.432689847121871; D=2024; E=0.1172571542027212},
    @{Row=19; A=45618; B=2024; C=2.033479419175133; D=2025; E=1.317145539573517}
)

# Copy the date number-format/border style from the last existing data row (A18)
# down onto the newly-created row 19, column A, before overwriting values so
# that the new row matches the established look of the table.
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

foreach ($r in $targetRows) {
    $ws.Cells.Item($r.Row, 1).Value() = $r.A
    $ws.Cells.Item($r.Row, 2).Value() = $r.B
    $ws.Cells.Item($r.Row, 3).Value() = $r.C
    $ws.Cells.Item($r.Row, 4).Value() = $r.D
    $ws.Cells.Item($r.Row, 5).Value() = $r.E
}
